# ============================================================================
# Module 3 (Strings) addition + small Module 2 wording tweaks.
# ============================================================================

$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------------
# 1. Slide 19 ("Passing Arrays" slide, sldId 275): give the empty title its
#    text back: "Passing Arrays " in Times New Roman.
# ----------------------------------------------------------------------------
$slide19 = $p.Slides.Item(19)
$title19 = $slide19.Shapes.Item("Title 1")
$titleRange = $title19.TextFrame.TextRange
$titleRange.Text = "Passing Arrays "
$titleRange.Font.Name = "Times New Roman"
$titleRange.Font.NameComplexScript = "Times New Roman"

# ----------------------------------------------------------------------------
# 2. Slide 20 (MODULE 2 cover slide, sldId 276): recolor the middle clause of
#    the subtitle paragraph green, splitting the trailing run into three.
# ----------------------------------------------------------------------------
$slide20 = $p.Slides.Item(20)
$subtitle20 = $slide20.Shapes.Item("Subtitle 2")
$subtitleRange = $subtitle20.TextFrame.TextRange
$fullSubtitle = $subtitleRange.Text
$needle = "Develop C programs to perform different operations on strings"
$idx = $fullSubtitle.IndexOf($needle)
$greenRun = $subtitleRange.Characters($idx + 1, $needle.Length)
$greenRun.Font.Color.RGB = 5287936

# ----------------------------------------------------------------------------
# 3. New slide 21 (sldId 277): "String" definition slide.
# ----------------------------------------------------------------------------
$slide21 = $p.Slides.Add(21, 2)

$title21 = $slide21.Shapes.Item(1)
$title21.TextFrame.TextRange.Text = "String`t"
$title21.TextFrame.TextRange.Font.Name = "Times New Roman"
$title21.TextFrame.TextRange.Font.NameComplexScript = "Times New Roman"
$title21.TextFrame.TextRange.Font.Bold = 1

$body21 = $slide21.Shapes.Item(2)
$body21Tf = $body21.TextFrame.TextRange
$body21Tf.Text = "String is a sequence of characters  stored in a `tcontiguous block of memory and terminated by a null character(" + [char]0x2018 + "\0" + [char]0x2019 + ")." + "`rThis null terminator is crucial as it signals the end of the string, allowing functions to know where the string ends.`rchar str_array[] = `"Hello`";"
$body21Tf.Font.Name = "Times New Roman"
$body21Tf.Font.NameComplexScript = "Times New Roman"

# Paragraph 1: recolor the "null character(...)." tail.
$p1Text = $body21Tf.Text
$needle1 = "null character("
$idx1 = $p1Text.IndexOf($needle1)
$tailLen = $p1Text.IndexOf("`r") - $idx1
$tailRun = $body21Tf.Characters($idx1 + 1, $tailLen)
$tailRun.Font.Color.ObjectThemeColor = 13
$tailRun.Font.Color.TintAndShade = 0.5

# Paragraph 3: "char str_array[] = "Hello";" -> split str_array out.
$p3Start = $p1Text.Length + ("This null terminator is crucial as it signals the end of the string, allowing functions to know where the string ends.").Length + 2
$fullText3 = $body21Tf.Text
$charIdx = $fullText3.IndexOf("str_array")
$strArrayRun = $body21Tf.Characters($charIdx + 1, "str_array".Length)
$strArrayRun.Font.Name = "Times New Roman"

# ----------------------------------------------------------------------------
# 4. New slide 22 (sldId 278): "String Operations" slide with sample code and
#    a floating "Sample Programs" label textbox.
# ----------------------------------------------------------------------------
$slide22 = $p.Slides.Add(22, 2)

$title22 = $slide22.Shapes.Item(1)
$title22Tf = $title22.TextFrame.TextRange
$title22Tf.Text = "String Operations`t"
$boldPart = $title22Tf.Characters(1, "String Operations".Length)
$boldPart.Font.Name = "Times New Roman"
$boldPart.Font.NameComplexScript = "Times New Roman"
$boldPart.Font.Bold = 1

$body22 = $slide22.Shapes.Item(2)
$body22Tf = $body22.TextFrame.TextRange
$body22Tf.Text = "String manipulation in C is done using functions from the <string.h> header file.`rConcatenation`rstrcat(str1, str2);`rCopying`rstrcpy(destination, source);`rLength Calculation`rint length = strlen(str);`rComparison`rint result1 = strcmp(str1, str2);`r"

$paras22 = $body22Tf.Paragraphs()
# Paragraph indices (1-based): 1 intro, 2 Concatenation, 3 strcat(...), 4 Copying,
# 5 strcpy(...), 6 Length Calculation, 7 int length = strlen(str);, 8 Comparison,
# 9 int result1 = strcmp(str1, str2);, 10 trailing blank paragraph.
$headingIdx = 2,4,6,8
foreach ($i in $headingIdx) {
  $para = $paras22.Item($i)
  $para.IndentLevel = 0
  $para.Font.Color.ObjectThemeColor = 13
  $para.Font.Color.TintAndShade = -0.1
}
$codeIdx = 3,5,7,9
foreach ($i in $codeIdx) {
  $para = $paras22.Item($i)
  $para.IndentLevel = 1
}

# New floating textbox labelling the sample-program group of slides.
$label = $slide22.Shapes.AddTextbox(1, 8917858/12700, 5034116/12700, 1938929/12700, 369332/12700)
$label.TextFrame.WordWrap = 0
$label.TextFrame.AutoSize = 1
$label.Fill.Visible = 0
$label.TextFrame.TextRange.Text = "Sample Programs"
